$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.366.78"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.626.62"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'597.42"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'168.47"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "2.627.68"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "'0.363"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "'5.25"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "'27.78"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "3.104.77"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "67.340.08"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "2.633.15"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "'12.16"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").Value = "'8.08"
$ws.Range("E20").Value = "  +7.82%  "
$ws.Range("D21").Value = "'356.85"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "'4.69"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "'10.52"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'1.93"
$ws.Range("E26").Value = "  -4.80%  "
$ws.Range("D27").Value = "'69.67"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "2.759.08"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "'0.0000101"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'551.20"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'7.95"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("D36").Value = "'0.996"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").Value = "'157.56"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").Value = "'19.01"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "'0.367"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.21"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'18.15"
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").Value = "'40.13"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "0.0₆0296"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'152.06"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'0.580"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'3.79"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "'1.71"
$ws.Range("E51").Value = "  -0.87%  "
